$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D23").Value = 84797
$ws.Range("E23").Value = "love is in the air  beef fondue   sauces"
$ws.Range("D24").Value = 109439
$ws.Range("E24").Value = "berry  good sandwich spread"
$ws.Range("D25").Value = 42522
$ws.Range("E25").Value = "the man s  taco dip"
$ws.Range("D26").Value = 62368
$ws.Range("E26").Value = "the best  chocolate chip cheesecake ever"
$ws.Range("D33").Value = 76808
$ws.Range("E33").Value = "cream  of spinach soup"
$ws.Range("D34").Value = 93249
$ws.Range("E34").Value = "grilled  ranch bread"
$ws.Range("D36").Value = 112959
$ws.Range("E36").Value = "sour cream  avocado dip  vegan"
$ws.Range("D166").Value = 38276
$ws.Range("E166").Value = "now and later  vegetarian empanadas"
$ws.Range("D174").Value = 62368
$ws.Range("E174").Value = "the best  chocolate chip cheesecake ever"
$ws.Range("D175").Value = 35964
$ws.Range("E175").Value = "rich  hot fudge cake"
$ws.Range("D176").Value = 24701
$ws.Range("E176").Value = "cream  of spinach soup  vegan"
$ws.Range("D208").Value = 83025
$ws.Range("E208").Value = "jeanne s style  birthday cake"
$ws.Range("D209").Value = 35964
$ws.Range("E209").Value = "rich  hot fudge cake"
$ws.Range("D210").Value = 32271
$ws.Range("E210").Value = "one pot  brownies"
$ws.Range("D211").Value = 58651
$ws.Range("E211").Value = "turtle  squares"
$ws.Range("D233").Value = 62368
$ws.Range("E233").Value = "the best  chocolate chip cheesecake ever"
$ws.Range("D234").Value = 39363
$ws.Range("E234").Value = "the best  banana bread  or muffins"
$ws.Range("D244").Value = 62368
$ws.Range("E244").Value = "the best  chocolate chip cheesecake ever"
$ws.Range("D245").Value = 71635
$ws.Range("E245").Value = "no bake  cookie crumble cheesecake"
$ws.Range("D246").Value = 49262
$ws.Range("E246").Value = "easiest ever  hollandaise sauce"
$ws.Range("D271").Value = 41756
$ws.Range("E271").Value = "souper  easy sweet   sour meatballs"
$ws.Range("D304").Value = 112140
$ws.Range("E304").Value = "all in the kitchen  chili"
$ws.Range("D314").Value = 93249
$ws.Range("E314").Value = "grilled  ranch bread"
$ws.Range("D315").Value = 58224
$ws.Range("E315").Value = "immoral  sandwich filling  loose meat"
$ws.Range("D316").Value = 99024
$ws.Range("E316").Value = "smoked  salmon  cracker spread"
$ws.Range("D319").Value = 44045
$ws.Range("E319").Value = "mennonite  corn fritters"
$ws.Range("D320").Value = 87098
$ws.Range("E320").Value = "homemade  vegetable soup from a can"
$ws.Range("D321").Value = 112140
$ws.Range("E321").Value = "all in the kitchen  chili"
$ws.Range("D348").Value = 87098
$ws.Range("E348").Value = "homemade  vegetable soup from a can"
$ws.Range("D349").Value = 83133
$ws.Range("E349").Value = "stove top  bbq  beef or pork ribs"
$ws.Range("D351").Value = 112140
$ws.Range("E351").Value = "all in the kitchen  chili"
$ws.Range("D354").Value = 23933
$ws.Range("E354").Value = "chinese  candy"
$ws.Range("D384").Value = 23933
$ws.Range("E384").Value = "chinese  candy"
$ws.Range("D388").Value = 112140
$ws.Range("E388").Value = "all in the kitchen  chili"
$ws.Range("D391").Value = 58224
$ws.Range("E391").Value = "immoral  sandwich filling  loose meat"
$ws.Range("D410").Value = 74805
$ws.Range("E410").Value = "never weep  whipped cream"
$ws.Range("D411").Value = 52804
$ws.Range("E411").Value = "jiffy  extra moist carrot cake"
$ws.Range("D417").Value = 63593
$ws.Range("E417").Value = "more  more    apple pear jigglers"
$ws.Range("D418").Value = 22123
$ws.Range("E418").Value = "i don t feel like cooking tonight  casserole"
$ws.Range("D419").Value = 112140
$ws.Range("E419").Value = "all in the kitchen  chili"
$ws.Range("D420").Value = 39959
$ws.Range("E420").Value = "calm your nerves  tonic"
$ws.Range("D421").Value = 58224
$ws.Range("E421").Value = "immoral  sandwich filling  loose meat"
$ws.Range("D443").Value = 75452
$ws.Range("E443").Value = "beat this  banana bread"
$ws.Range("D444").Value = 83062
$ws.Range("E444").Value = "spicy  banana bread"
$ws.Range("D445").Value = 39363
$ws.Range("E445").Value = "the best  banana bread  or muffins"
$ws.Range("D446").Value = 95926
$ws.Range("E446").Value = "say what   banana sandwich"
$ws.Range("D550").Value = 59952
$ws.Range("E550").Value = "global gourmet  taco casserole"
$ws.Range("D551").Value = 44123
$ws.Range("E551").Value = "george s at the cove  black bean soup"
$ws.Range("D618").Value = 38276
$ws.Range("E618").Value = "now and later  vegetarian empanadas"
$ws.Range("D620").Value = 67888
$ws.Range("E620").Value = "backyard style  barbecued ribs"
$ws.Range("D621").Value = 64045
$ws.Range("E621").Value = "some like it hot"
$ws.Range("D653").Value = 42570
$ws.Range("E653").Value = "pick me up  party chicken kabobs"
$ws.Range("D654").Value = 58224
$ws.Range("E654").Value = "immoral  sandwich filling  loose meat"
$ws.Range("D769").Value = 63793
$ws.Range("E769").Value = "tide me over   indian chaat  simple veggie salad"
$ws.Range("D770").Value = 95926
$ws.Range("E770").Value = "say what   banana sandwich"
$ws.Range("D775").Value = 30131
$ws.Range("E775").Value = "momma s special  marinade"
$ws.Range("D776").Value = 93249
$ws.Range("E776").Value = "grilled  ranch bread"
$ws.Range("D803").Value = 112959
$ws.Range("E803").Value = "sour cream  avocado dip  vegan"
$ws.Range("D804").Value = 111875
$ws.Range("E804").Value = "the elvis  smoothie"
$ws.Range("D805").Value = 23850
$ws.Range("E805").Value = "cream  of cauliflower soup  vegan"
$ws.Range("D828").Value = 83133
$ws.Range("E828").Value = "stove top  bbq  beef or pork ribs"
$ws.Range("D829").Value = 64302
$ws.Range("E829").Value = "red  macaroni salad"
$ws.Range("D875").Value = 59952
$ws.Range("E875").Value = "global gourmet  taco casserole"
$ws.Range("D876").Value = 32169
$ws.Range("E876").Value = "make that chicken dance  salsa pasta"
$ws.Range("D907").Value = 23850
$ws.Range("E907").Value = "cream  of cauliflower soup  vegan"
$ws.Range("D908").Value = 137739
$ws.Range("E908").Value = "arriba   baked winter squash mexican style"
$ws.Range("D909").Value = 59534
$ws.Range("E909").Value = "twisted american chop suey"
$ws.Range("D942").Value = 83873
$ws.Range("E942").Value = "crispy crunchy  chicken"
$ws.Range("D943").Value = 107229
$ws.Range("E943").Value = "open sesame  noodles"
$ws.Range("D944").Value = 63986
$ws.Range("E944").Value = "chicken lickin  good  pork chops"
$ws.Range("D945").Value = 98930
$ws.Range("E945").Value = "steamed  chicken cutlets in packages"
$ws.Range("D946").Value = 54272
$ws.Range("E946").Value = "fool the meat eaters  chili"
$ws.Range("D968").Value = 32169
$ws.Range("E968").Value = "make that chicken dance  salsa pasta"
$ws.Range("D969").Value = 53402
$ws.Range("E969").Value = "killer  lasagna"
$ws.Range("D970").Value = 94710
$ws.Range("E970").Value = "italian  fries"
$ws.Range("D971").Value = 47366
$ws.Range("E971").Value = "forgotten  minestrone"
$ws.Range("D1049").Value = 30300
$ws.Range("E1049").Value = "munch without guilt  tomatoes"
$ws.Range("D1051").Value = 95926
$ws.Range("E1051").Value = "say what   banana sandwich"
$ws.Range("D1073").Value = 59534
$ws.Range("E1073").Value = "twisted american chop suey"
$ws.Range("D1074").Value = 41756
$ws.Range("E1074").Value = "souper  easy sweet   sour meatballs"
$ws.Range("D1075").Value = 112140
$ws.Range("E1075").Value = "all in the kitchen  chili"
$ws.Range("D1086").Value = 67664
$ws.Range("E1086").Value = "healthy for them  yogurt popsicles"
$ws.Range("D1168").Value = 52804
$ws.Range("E1168").Value = "jiffy  extra moist carrot cake"
$ws.Range("D1169").Value = 26995
$ws.Range("E1169").Value = "keep it going  german friendship cake"
$ws.Range("D1171").Value = 27087
$ws.Range("E1171").Value = "get the sensation  brownies"
$ws.Range("D1180").Value = 26835
$ws.Range("E1180").Value = "one bowl  perfect pound cake"
$ws.Range("D1181").Value = 75452
$ws.Range("E1181").Value = "beat this  banana bread"
$ws.Range("D1223").Value = 24701
$ws.Range("E1223").Value = "cream  of spinach soup  vegan"
$ws.Range("D1224").Value = 25274
$ws.Range("E1224").Value = "aww  marinated olives"
$ws.Range("D1225").Value = 54272
$ws.Range("E1225").Value = "fool the meat eaters  chili"
$ws.Range("D1226").Value = 81185
$ws.Range("E1226").Value = "mock a mole   low fat guacamole"
$ws.Range("D1261").Value = 107699
$ws.Range("E1261").Value = "deep fried dessert thingys"
$ws.Range("D1319").Value = 62368
$ws.Range("E1319").Value = "the best  chocolate chip cheesecake ever"
$ws.Range("D1320").Value = 38276
$ws.Range("E1320").Value = "now and later  vegetarian empanadas"
$ws.Range("D1321").Value = 35653
$ws.Range("E1321").Value = "make it your way  shortcakes"
$ws.Range("D1363").Value = 60219
$ws.Range("E1363").Value = "mexican pasta"
$ws.Range("D1364").Value = 25274
$ws.Range("E1364").Value = "aww  marinated olives"
$ws.Range("D1365").Value = 32169
$ws.Range("E1365").Value = "make that chicken dance  salsa pasta"
$ws.Range("D1366").Value = 54100
$ws.Range("E1366").Value = "grilled  venison burgers"
$ws.Range("D1468").Value = 31490
$ws.Range("E1468").Value = "a bit different  breakfast pizza"
$ws.Range("D1469").Value = 103948
$ws.Range("E1469").Value = "smells like sunday  chicken fricassee with meatballs"
$ws.Range("D1470").Value = 25775
$ws.Range("E1470").Value = "how i got my family to eat spinach  spinach casserole"
$ws.Range("D1471").Value = 39363
$ws.Range("E1471").Value = "the best  banana bread  or muffins"
$ws.Range("D1504").Value = 112959
$ws.Range("E1504").Value = "sour cream  avocado dip  vegan"
$ws.Range("D1505").Value = 27087
$ws.Range("E1505").Value = "get the sensation  brownies"
